# Add two new IOP DMAC channel "interrupt-style" registers (STAT @ row159,
# MASK @ row171) each with 11 bitfields: VBLNK,GPU,CDROM,DMA,TMR0,TMR1,TMR2,
# CON_MC,SIO,SPU,PIO - mirroring the existing BCR/CHCR register blocks above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Field name / bit-position table shared by both new registers.
$fields = @(
    @{ Name = "VBLNK";  Pos = 0 },
    @{ Name = "GPU";    Pos = 1 },
    @{ Name = "CDROM";  Pos = 2 },
    @{ Name = "DMA";    Pos = 3 },
    @{ Name = "TMR0";   Pos = 4 },
    @{ Name = "TMR1";   Pos = 5 },
    @{ Name = "TMR2";   Pos = 6 },
    @{ Name = "CON_MC"; Pos = 7 },
    @{ Name = "SIO";    Pos = 8 },
    @{ Name = "SPU";    Pos = 9 },
    @{ Name = "PIO";    Pos = 10 }
)

function Add-IOPDmacRegister($registerName, $startRow) {
    $endRow = $startRow + $fields.Length - 1

    # Column A: register name only on the first row of the block.
    $ws.Range("A$startRow").Value = $registerName

    # Columns B-F: field index / name / startpos / length / initial value.
    for ($i = 0; $i -lt $fields.Length; $i++) {
        $r = $startRow + $i
        $ws.Range("B$r").Value = $fields[$i].Pos
        $ws.Range("C$r").Value = $fields[$i].Name
        $ws.Range("D$r").Value = $fields[$i].Pos
        $ws.Range("E$r").Value = 1
        $ws.Range("F$r").Value = 0
    }

    # Column H: register forward declaration (single cell formula).
    $ws.Range("H$startRow").Formula = '="class IOPDmacChannelRegister_"&A' + $startRow + '&"_t;"'

    # Column I: field declaration (shared formula across the whole block).
    $ws.Range("I$startRow`:I$endRow").Formula = '="static constexpr u8 "&C' + $startRow + '&" = "&B' + $startRow + '&";"'

    # Column J: field init / registerField call (shared formula across the whole block).
    $ws.Range("J$startRow`:J$endRow").Formula = '="registerField(Fields::"&C' + $startRow + '&", """&C' + $startRow + '&""", "&D' + $startRow + '&", "&E' + $startRow + '&", "&F' + $startRow + '&");"'
}

# New row 159-169 block: STAT register.
Add-IOPDmacRegister "STAT" 159

# New row 171-181 block: MASK register (row 170 left blank, matching the
# existing convention of a blank separator row between register blocks).
Add-IOPDmacRegister "MASK" 171

# Match the author's final selection/view state.
$ws.Range("B171:J182").Select()
